$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the existing single-column (B) checklist over to column C ---
# (Excel "cut & paste one column to the right" style move)
$ws.Range("C1").Value = "Contract PDF on ContractPage"
$ws.Range("C2").Value = "Email Contract"
$ws.Range("C3").Value = "Save Contract as PDF to Membership Accounting"
$ws.Range("C4").Value = "Addons per Club in table"
$ws.Range("C5").Value = "Online Windows Server"
$ws.Range("C6").Value = "Connect to Existing Wellbridge Website"
$ws.Range("C7").Value = "Credit card process"
$ws.Range("C8").Value = "Update text on Enrollment Confirmation"
$ws.Range("C9").Value = "Ads for PT?"
$ws.Range("C10").Value = "Credit Card Test Account"

# Clear out the old column B (all of its content has now moved to C / K)
$ws.Range("B1:B14").ClearContents()

# --- New BOSS / transaction related checklist items (typed in this order) ---
$ws.Range("C11").Value = "Submit Transaction Data to BOSS"
$ws.Range("D12").Value = "Create asptheade transaction number"
$ws.Range("D13").Value = "enter asptitemd and aspttendd after"
$ws.Range("C14").Value = "Enter Token into strcustr"
$ws.Range("C15").Value = "Save CC last 4 digits and cvv exp date"
$ws.Range("C16").Value = "Create ONLINE Boss account"
$ws.Range("C17").Value = "Script to OPEN / CLOSE Drawer"
$ws.Range("C18").Value = "Ifee for Online? JoAnna? Add a For Web column in ifee table?"

# --- Row 16: mark this one done with an "X" and today's completion date ---
$ws.Range("A16").Value = "X"
$ws.Range("B16").Value = 45778
$ws.Range("B16").NumberFormat = "m/d/yyyy"

# --- Cleanup list moved into columns J/K ---
$ws.Range("J1").Value = "CLEANUP"
$ws.Range("K2").Value = "Remove extra text from shopping cart"
$ws.Range("K3").Value = "Remove CC Processor text from PaymentPage"
$ws.Range("K4").Value = "Remove any logs that may show data in the browser"
$ws.Range("K5").Value = "Remove Payment Method CC on Contract Page before contract"

# --- Selection matches the author's final cursor position ---
$ws.Range("J11").Select()
